$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("D1").Value = "Average Memory Usage (bytes)"

# Update Average Run Time (ms) column C values
$ws.Range("C2").Value = 17.599
$ws.Range("C3").Value = 16.262
$ws.Range("C4").Value = 18.036
$ws.Range("C5").Value = 16.223
$ws.Range("C6").Value = 18.094
$ws.Range("C7").Value = 18.486
$ws.Range("C8").Value = 18.107
$ws.Range("C9").Value = 16.079
$ws.Range("C10").Value = 16.241
$ws.Range("C11").Value = 16.622
$ws.Range("C12").Value = 16.828
$ws.Range("C13").Value = 16.807
$ws.Range("C14").Value = 17.333
$ws.Range("C15").Value = 16.832
$ws.Range("C16").Value = 17.522

# Update Average Memory Usage column D value for row 2 (others unchanged)
$ws.Range("D2").Value = 5608.2
